# Weekly update: a new price observation is inserted at row 24 (pushing the
# existing rows 24-40 down to 25-41) for
# "Hortaliza, Comercializadora del Agro de Limarí - Zapallo italiano".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24; Excel shifts rows 24:40 down to 25:41 and
# inherits the formatting (incl. the date style on column D) from the row
# above, same as the native "Insert" ribbon action.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new observation.
$ws.Range("A24").Value = 2
$ws.Range("B24").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44546
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 100112032
$ws.Range("G24").Value = "Zapallo italiano"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 8000
$ws.Range("M24").Value = 7500
$ws.Range("N24").Value = "`$/caja 60 unidades"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 125
$ws.Range("Q24").Value = 60
$ws.Range("R24").Value = "Hortaliza"
